$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Add the new row of data (row 10)
# NOTE: write B10 before A10 so the shared-strings table gets the new
# unique strings in the same order as the target file (task text = 18, date = 19)
$ws.Range("B10").Value = "Fixed erros and reseached&tested better ways of applying Material-ui styles to custom components"
$ws.Range("A10").Value = "01.07.2019 - `n04.07.2019"
$ws.Range("C10").Value = 6

# Match formatting used by the other wrapped-text rows (style index 6: wrapText only)
$ws.Range("A10:B10").WrapText = $true
$ws.Rows.Item(10).AutoFit()

# Update view / selection state
$ws.Application.ActiveWindow.ScrollRow = 5
$ws.Range("A10").Select()
